# Replace the sole body paragraph ("TEST 3") with the full assignment-brief
# paragraph text, keeping a bold/18pt "TEST" run in the middle (where the
# original word "TEST" was) and the proofing marks (gramStart/gramEnd,
# spellStart/spellEnd) that Word leaves around phrases it flagged.
# InsertXML replaces the whole target Range's contents with the given
# WordprocessingML, so we pass the complete replacement <w:p> here.
$d = $word.ActiveDocument
$r = $d.Paragraphs(1).Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Ideal Job Find an advertisement showing what you believe to be your ideal job. This may require several years of experience, and hence be something that you must work towards, rather than something that you are ready for now, or will be able to fill as soon as you graduate. There are various ways to search for IT jobs, including websites like {seek.com.au}. You should include the following information. • The job advertisement itself. Include a link, and a snapshot of it (in case the link expires before the assignment deadline). • A description (in your own words) of the position, and particularly what makes this position appealing to you. • A description (in your own words) of the skills, qualifications and experience required for the position. • A description (in your own words) of the skills, qualifications and experience you currently have. </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>TEST</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>• A plan describing how you will obtain the skills, qualifications and experience required for the position, building on those you have now. This need not be greatly detailed, (and will probably change significantly over time anyway</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>), but</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> try to be as specific as you can. Personal Profile There are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a number of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> online tests that are commonly used by employers to get specific information about potential employees. One of the best-known of these is the Myers-Briggs Type Indicator (MBTI) test, which was developed by the mother-and-daughter team of Isabel Myers and Katharine Briggs. Another popular one (and very relevant to students) is about learning styles. There are various other tests available online as well, and while there is no guarantee that any specific test will be necessarily one that you will encounter in your career, it seems highly likely that you will have to do </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>some kind of test</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> like this as part of a recruitment process. Accordingly, you are required to present the following information: • The results of an online Myers-Briggs test. www.16personalities.com • The results of an online learning style test. • The results of one further online test of your choosing. The third test should be distinct from both the Myers-Briggs test and the learning styles </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tes</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$r.InsertXML($xml)
